$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 50.41084233333333
$ws.Cells.Item(2, 8).Value = 151.232527
$ws.Cells.Item(2, 9).Value = 0.1533822412306554
$ws.Cells.Item(2, 10).Value = 0.1625332948695183
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 14.03654733333333
$ws.Cells.Item(2, 14).Value = 42.10964199999999
$ws.Cells.Item(2, 15).Value = 0.3033248635033713
$ws.Cells.Item(2, 16).Value = 0.3241347232362796
$ws.Cells.Item(2, 17).Value = 707.594174525037
$ws.Cells.Item(2, 18).Value = 6368.347570725333
$ws.Cells.Item(2, 19).Value = 0.0465246473851297
$ws.Cells.Item(2, 20).Value = 0.05268268454921195

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 50.41084233333333
$ws.Cells.Item(3, 8).Value = 151.232527
$ws.Cells.Item(3, 9).Value = 0.1533822412306554
$ws.Cells.Item(3, 10).Value = 0.1625332948695183
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 5.85326
$ws.Cells.Item(3, 14).Value = 17.55978
$ws.Cells.Item(3, 15).Value = 0.1264868951307928
$ws.Cells.Item(3, 16).Value = 0.1351646359375356
$ws.Cells.Item(3, 17).Value = 295.0677669960066
$ws.Cells.Item(3, 18).Value = 2655.60990296406
$ws.Cells.Item(3, 19).Value = 0.01940084346146788
$ws.Cells.Item(3, 20).Value = 0.02196875362876657

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 50.41084233333333
$ws.Cells.Item(4, 8).Value = 151.232527
$ws.Cells.Item(4, 9).Value = 0.1533822412306554
$ws.Cells.Item(4, 10).Value = 0.1625332948695183
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 9.056607
$ws.Cells.Item(4, 14).Value = 27.169821
$ws.Cells.Item(4, 15).Value = 0.1957100999869824
$ws.Cells.Item(4, 16).Value = 0.2091369575218488
$ws.Cells.Item(4, 17).Value = 456.551187551963
$ws.Cells.Item(4, 18).Value = 4108.960687967667
$ws.Cells.Item(4, 19).Value = 0.03001845376747901
$ws.Cells.Item(4, 20).Value = 0.03399171878501258

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 50.41084233333333
$ws.Cells.Item(5, 8).Value = 151.232527
$ws.Cells.Item(5, 9).Value = 0.1533822412306554
$ws.Cells.Item(5, 10).Value = 0.1625332948695183
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 8.416348666666666
$ws.Cells.Item(5, 14).Value = 25.249046
$ws.Cells.Item(5, 15).Value = 0.1818743420221987
$ws.Cells.Item(5, 16).Value = 0.1943519856376384
$ws.Cells.Item(5, 17).Value = 424.2752256576935
$ws.Cells.Item(5, 18).Value = 3818.477030919242
$ws.Cells.Item(5, 19).Value = 0.0278962942017156
$ws.Cells.Item(5, 20).Value = 0.03158866859011868

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 50.41084233333333
$ws.Cells.Item(6, 8).Value = 151.232527
$ws.Cells.Item(6, 9).Value = 0.1533822412306554
$ws.Cells.Item(6, 10).Value = 0.1625332948695183
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 13).Value = 8.912860999999999
$ws.Cells.Item(6, 14).Value = 17.825722
$ws.Cells.Item(6, 15).Value = 0.1926037993566548
$ws.Cells.Item(6, 16).Value = 0.1372116976666974
$ws.Cells.Item(6, 17).Value = 449.3048306099157
$ws.Cells.Item(6, 18).Value = 2695.828983659494
$ws.Cells.Item(6, 19).Value = 0.02954200241486318
$ws.Cells.Item(6, 20).Value = 0.02230146931640853

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 72.399292
$ws.Cells.Item(7, 8).Value = 217.197876
$ws.Cells.Item(7, 9).Value = 0.2202852631789851
$ws.Cells.Item(7, 10).Value = 0.2334278684964451
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 14.03654733333333
$ws.Cells.Item(7, 14).Value = 42.10964199999999
$ws.Cells.Item(7, 15).Value = 0.3033248635033713
$ws.Cells.Item(7, 16).Value = 0.3241347232362796
$ws.Cells.Item(7, 17).Value = 1016.236089057821
$ws.Cells.Item(7, 18).Value = 9146.124801520391
$ws.Cells.Item(7, 19).Value = 0.06681799738556987
$ws.Cells.Item(7, 20).Value = 0.0756620775507299

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 72.399292
$ws.Cells.Item(8, 8).Value = 217.197876
$ws.Cells.Item(8, 9).Value = 0.2202852631789851
$ws.Cells.Item(8, 10).Value = 0.2334278684964451
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 5.85326
$ws.Cells.Item(8, 14).Value = 17.55978
$ws.Cells.Item(8, 15).Value = 0.1264868951307928
$ws.Cells.Item(8, 16).Value = 0.1351646359375356
$ws.Cells.Item(8, 17).Value = 423.77187989192
$ws.Cells.Item(8, 18).Value = 3813.94691902728
$ws.Cells.Item(8, 19).Value = 0.02786319898257939
$ws.Cells.Item(8, 20).Value = 0.03155119286299694

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 72.399292
$ws.Cells.Item(9, 8).Value = 217.197876
$ws.Cells.Item(9, 9).Value = 0.2202852631789851
$ws.Cells.Item(9, 10).Value = 0.2334278684964451
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 9.056607
$ws.Cells.Item(9, 14).Value = 27.169821
$ws.Cells.Item(9, 15).Value = 0.1957100999869824
$ws.Cells.Item(9, 16).Value = 0.2091369575218488
$ws.Cells.Item(9, 17).Value = 655.691934722244
$ws.Cells.Item(9, 18).Value = 5901.227412500196
$ws.Cells.Item(9, 19).Value = 0.04311205088241789
$ws.Cells.Item(9, 20).Value = 0.04881839421815674

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 72.399292
$ws.Cells.Item(10, 8).Value = 217.197876
$ws.Cells.Item(10, 9).Value = 0.2202852631789851
$ws.Cells.Item(10, 10).Value = 0.2334278684964451
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 8.416348666666666
$ws.Cells.Item(10, 14).Value = 25.249046
$ws.Cells.Item(10, 15).Value = 0.1818743420221987
$ws.Cells.Item(10, 16).Value = 0.1943519856376384
$ws.Cells.Item(10, 17).Value = 609.3376846918106
$ws.Cells.Item(10, 18).Value = 5484.039162226296
$ws.Cells.Item(10, 19).Value = 0.04006423729786479
$ws.Cells.Item(10, 20).Value = 0.04536716974544564

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 72.399292
$ws.Cells.Item(11, 8).Value = 217.197876
$ws.Cells.Item(11, 9).Value = 0.2202852631789851
$ws.Cells.Item(11, 10).Value = 0.2334278684964451
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 13).Value = 8.912860999999999
$ws.Cells.Item(11, 14).Value = 17.825722
$ws.Cells.Item(11, 15).Value = 0.1926037993566548
$ws.Cells.Item(11, 16).Value = 0.1372116976666974
$ws.Cells.Item(11, 17).Value = 645.284826094412
$ws.Cells.Item(11, 18).Value = 3871.708956566472
$ws.Cells.Item(11, 19).Value = 0.04242777863055315
$ws.Cells.Item(11, 20).Value = 0.03202903411911582

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 63.396933
$ws.Cells.Item(12, 8).Value = 190.190799
$ws.Cells.Item(12, 9).Value = 0.192894290605017
$ws.Cells.Item(12, 10).Value = 0.2044027024380561
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 14.03654733333333
$ws.Cells.Item(12, 14).Value = 42.10964199999999
$ws.Cells.Item(12, 15).Value = 0.3033248635033713
$ws.Cells.Item(12, 16).Value = 0.3241347232362796
$ws.Cells.Item(12, 17).Value = 889.8740508426619
$ws.Cells.Item(12, 18).Value = 8008.866457583957
$ws.Cells.Item(12, 19).Value = 0.05850963436834643
$ws.Cells.Item(12, 20).Value = 0.06625401338350695

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 63.396933
$ws.Cells.Item(13, 8).Value = 190.190799
$ws.Cells.Item(13, 9).Value = 0.192894290605017
$ws.Cells.Item(13, 10).Value = 0.2044027024380561
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 5.85326
$ws.Cells.Item(13, 14).Value = 17.55978
$ws.Cells.Item(13, 15).Value = 0.1264868951307928
$ws.Cells.Item(13, 16).Value = 0.1351646359375356
$ws.Cells.Item(13, 17).Value = 371.0787320515799
$ws.Cells.Item(13, 18).Value = 3339.70858846422
$ws.Cells.Item(13, 19).Value = 0.02439859990708547
$ws.Cells.Item(13, 20).Value = 0.02762801685968828

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 63.396933
$ws.Cells.Item(14, 8).Value = 190.190799
$ws.Cells.Item(14, 9).Value = 0.192894290605017
$ws.Cells.Item(14, 10).Value = 0.2044027024380561
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 9.056607
$ws.Cells.Item(14, 14).Value = 27.169821
$ws.Cells.Item(14, 15).Value = 0.1957100999869824
$ws.Cells.Item(14, 16).Value = 0.2091369575218488
$ws.Cells.Item(14, 17).Value = 574.1611071863309
$ws.Cells.Item(14, 18).Value = 5167.449964676979
$ws.Cells.Item(14, 19).Value = 0.03775136090122592
$ws.Cells.Item(14, 20).Value = 0.04274815929713884

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 63.396933
$ws.Cells.Item(15, 8).Value = 190.190799
$ws.Cells.Item(15, 9).Value = 0.192894290605017
$ws.Cells.Item(15, 10).Value = 0.2044027024380561
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 8.416348666666666
$ws.Cells.Item(15, 14).Value = 25.249046
$ws.Cells.Item(15, 15).Value = 0.1818743420221987
$ws.Cells.Item(15, 16).Value = 0.1943519856376384
$ws.Cells.Item(15, 17).Value = 533.570692525306
$ws.Cells.Item(15, 18).Value = 4802.136232727754
$ws.Cells.Item(15, 19).Value = 0.03508252218362626
$ws.Cells.Item(15, 20).Value = 0.03972607108853556

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 63.396933
$ws.Cells.Item(16, 8).Value = 190.190799
$ws.Cells.Item(16, 9).Value = 0.192894290605017
$ws.Cells.Item(16, 10).Value = 0.2044027024380561
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 13).Value = 8.912860999999999
$ws.Cells.Item(16, 14).Value = 17.825722
$ws.Cells.Item(16, 15).Value = 0.1926037993566548
$ws.Cells.Item(16, 16).Value = 0.1372116976666974
$ws.Cells.Item(16, 17).Value = 565.048051655313
$ws.Cells.Item(16, 18).Value = 3390.288309931878
$ws.Cells.Item(16, 19).Value = 0.03715217324473297
$ws.Cells.Item(16, 20).Value = 0.02804644180918647

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 86.94092833333333
$ws.Cells.Item(17, 8).Value = 260.822785
$ws.Cells.Item(17, 9).Value = 0.2645302840659494
$ws.Cells.Item(17, 10).Value = 0.2803126249625782
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 14.03654733333333
$ws.Cells.Item(17, 14).Value = 42.10964199999999
$ws.Cells.Item(17, 15).Value = 0.3033248635033713
$ws.Cells.Item(17, 16).Value = 0.3241347232362796
$ws.Cells.Item(17, 17).Value = 1220.350455754774
$ws.Cells.Item(17, 18).Value = 10983.15410179297
$ws.Cells.Item(17, 19).Value = 0.08023861230681212
$ws.Cells.Item(17, 20).Value = 0.09085905511188033

$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 7).Value = 86.94092833333333
$ws.Cells.Item(18, 8).Value = 260.822785
$ws.Cells.Item(18, 9).Value = 0.2645302840659494
$ws.Cells.Item(18, 10).Value = 0.2803126249625782
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 13).Value = 5.85326
$ws.Cells.Item(18, 14).Value = 17.55978
$ws.Cells.Item(18, 15).Value = 0.1264868951307928
$ws.Cells.Item(18, 16).Value = 0.1351646359375356
$ws.Cells.Item(18, 17).Value = 508.8878581763666
$ws.Cells.Item(18, 18).Value = 4579.9907235873
$ws.Cells.Item(18, 19).Value = 0.03345961429956858
$ws.Cells.Item(18, 20).Value = 0.03788835390176184

$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 7).Value = 86.94092833333333
$ws.Cells.Item(19, 8).Value = 260.822785
$ws.Cells.Item(19, 9).Value = 0.2645302840659494
$ws.Cells.Item(19, 10).Value = 0.2803126249625782
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 13).Value = 9.056607
$ws.Cells.Item(19, 14).Value = 27.169821
$ws.Cells.Item(19, 15).Value = 0.1957100999869824
$ws.Cells.Item(19, 16).Value = 0.2091369575218488
$ws.Cells.Item(19, 17).Value = 787.3898201301649
$ws.Cells.Item(19, 18).Value = 7086.508381171485
$ws.Cells.Item(19, 19).Value = 0.0517712483441318
$ws.Cells.Item(19, 20).Value = 0.05862372953963664

$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 7).Value = 86.94092833333333
$ws.Cells.Item(20, 8).Value = 260.822785
$ws.Cells.Item(20, 9).Value = 0.2645302840659494
$ws.Cells.Item(20, 10).Value = 0.2803126249625782
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 13).Value = 8.416348666666666
$ws.Cells.Item(20, 14).Value = 25.249046
$ws.Cells.Item(20, 15).Value = 0.1818743420221987
$ws.Cells.Item(20, 16).Value = 0.1943519856376384
$ws.Cells.Item(20, 17).Value = 731.7251662570121
$ws.Cells.Item(20, 18).Value = 6585.52649631311
$ws.Cells.Item(20, 19).Value = 0.04811127135943986
$ws.Cells.Item(20, 20).Value = 0.05447931526077571

$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 7).Value = 86.94092833333333
$ws.Cells.Item(21, 8).Value = 260.822785
$ws.Cells.Item(21, 9).Value = 0.2645302840659494
$ws.Cells.Item(21, 10).Value = 0.2803126249625782
$ws.Cells.Item(21, 11).Value = 2
$ws.Cells.Item(21, 13).Value = 8.912860999999999
$ws.Cells.Item(21, 14).Value = 17.825722
$ws.Cells.Item(21, 15).Value = 0.1926037993566548
$ws.Cells.Item(21, 16).Value = 0.1372116976666974
$ws.Cells.Item(21, 17).Value = 774.8924094459616
$ws.Cells.Item(21, 18).Value = 4649.35445667577
$ws.Cells.Item(21, 19).Value = 0.05094953775599702
$ws.Cells.Item(21, 20).Value = 0.03846217114852361

$ws.Cells.Item(22, 5).Value = 2
$ws.Cells.Item(22, 7).Value = 55.513536
$ws.Cells.Item(22, 8).Value = 111.027072
$ws.Cells.Item(22, 9).Value = 0.1689079209193933
$ws.Cells.Item(22, 10).Value = 0.1193235092334022
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 13).Value = 14.03654733333333
$ws.Cells.Item(22, 14).Value = 42.10964199999999
$ws.Cells.Item(22, 15).Value = 0.3033248635033713
$ws.Cells.Item(22, 16).Value = 0.3241347232362796
$ws.Cells.Item(22, 17).Value = 779.2183757047039
$ws.Cells.Item(22, 18).Value = 4675.310254228223
$ws.Cells.Item(22, 19).Value = 0.05123397205751321
$ws.Cells.Item(22, 20).Value = 0.03867689264095046

$ws.Cells.Item(23, 5).Value = 2
$ws.Cells.Item(23, 7).Value = 55.513536
$ws.Cells.Item(23, 8).Value = 111.027072
$ws.Cells.Item(23, 9).Value = 0.1689079209193933
$ws.Cells.Item(23, 10).Value = 0.1193235092334022
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 13).Value = 5.85326
$ws.Cells.Item(23, 14).Value = 17.55978
$ws.Cells.Item(23, 15).Value = 0.1264868951307928
$ws.Cells.Item(23, 16).Value = 0.1351646359375356
$ws.Cells.Item(23, 17).Value = 324.93515972736
$ws.Cells.Item(23, 18).Value = 1949.61095836416
$ws.Cells.Item(23, 19).Value = 0.02136463848009155
$ws.Cells.Item(23, 20).Value = 0.01612831868432197

$ws.Cells.Item(24, 5).Value = 2
$ws.Cells.Item(24, 7).Value = 55.513536
$ws.Cells.Item(24, 8).Value = 111.027072
$ws.Cells.Item(24, 9).Value = 0.1689079209193933
$ws.Cells.Item(24, 10).Value = 0.1193235092334022
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 13).Value = 9.056607
$ws.Cells.Item(24, 14).Value = 27.169821
$ws.Cells.Item(24, 15).Value = 0.1957100999869824
$ws.Cells.Item(24, 16).Value = 0.2091369575218488
$ws.Cells.Item(24, 17).Value = 502.764278732352
$ws.Cells.Item(24, 18).Value = 3016.585672394112
$ws.Cells.Item(24, 19).Value = 0.03305698609172777
$ws.Cells.Item(24, 20).Value = 0.02495495568190396

$ws.Cells.Item(25, 5).Value = 2
$ws.Cells.Item(25, 7).Value = 55.513536
$ws.Cells.Item(25, 8).Value = 111.027072
$ws.Cells.Item(25, 9).Value = 0.1689079209193933
$ws.Cells.Item(25, 10).Value = 0.1193235092334022
$ws.Cells.Item(25, 11).Value = 3
$ws.Cells.Item(25, 13).Value = 8.416348666666666
$ws.Cells.Item(25, 14).Value = 25.249046
$ws.Cells.Item(25, 15).Value = 0.1818743420221987
$ws.Cells.Item(25, 16).Value = 0.1943519856376384
$ws.Cells.Item(25, 17).Value = 467.221274695552
$ws.Cells.Item(25, 18).Value = 2803.327648173312
$ws.Cells.Item(25, 19).Value = 0.03072001697955223
$ws.Cells.Item(25, 20).Value = 0.02319076095276279

$ws.Cells.Item(26, 5).Value = 2
$ws.Cells.Item(26, 7).Value = 55.513536
$ws.Cells.Item(26, 8).Value = 111.027072
$ws.Cells.Item(26, 9).Value = 0.1689079209193933
$ws.Cells.Item(26, 10).Value = 0.1193235092334022
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 13).Value = 8.912860999999999
$ws.Cells.Item(26, 14).Value = 17.825722
$ws.Cells.Item(26, 15).Value = 0.1926037993566548
$ws.Cells.Item(26, 16).Value = 0.1372116976666974
$ws.Cells.Item(26, 17).Value = 494.784429986496
$ws.Cells.Item(26, 18).Value = 1979.137719945984
$ws.Cells.Item(26, 19).Value = 0.03253230731050855
$ws.Cells.Item(26, 20).Value = 0.01637258127346295

